$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

$ws.Range("A11").Value = "T010"
$ws.Range("B11").Value = "E005"
$ws.Range("C11").Value = "A004"
$ws.Range("D11").Value = "dhj"
$ws.Range("E11").Value = "cvbn"
$ws.Range("F11").Value = "bdm"
$ws.Range("G11").Value = "cdccdcc"
$ws.Range("H11").Value = "eeee"
$ws.Range("I11").Value = "efcrcrdv"
$ws.Range("J11").Value = "ev"
